# Atualizado por script em 29-11-2023 20:45
# Adds the new match row (row 82) to the Armenia Premier League 2023-2024 sheet,
# mirroring the formatting of the previous data row and filling in the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing data row (81) into the new row (82)
$ws.Range("A81:V81").Copy()
$ws.Range("A82:V82").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row with the match data
$ws.Range("A82").Value = 81
$ws.Range("B82").Value = "armenia"
$ws.Range("C82").Value = "premier-league"
$ws.Range("D82").Value = "2023-2024"
$ws.Range("E82").Value = 45259.625
$ws.Range("F82").Value = "Urartu"
$ws.Range("G82").Value = 2
$ws.Range("H82").Value = "Ararat Yerevan"
$ws.Range("I82").Value = 1
$ws.Range("J82").Value = 1.41
$ws.Range("K82").Value = "28/11/2023 03:12"
$ws.Range("L82").Value = 1.49
$ws.Range("M82").Value = "29/11/2023 14:59"
$ws.Range("N82").Value = 4.41
$ws.Range("O82").Value = "28/11/2023 03:12"
$ws.Range("P82").Value = 4.47
$ws.Range("Q82").Value = "29/11/2023 14:59"
$ws.Range("R82").Value = 6.04
$ws.Range("S82").Value = "28/11/2023 03:12"
$ws.Range("T82").Value = 6.18
$ws.Range("U82").Value = "29/11/2023 14:59"
$ws.Range("V82").Value = "https://www.betexplorer.com/football/armenia/premier-league/urartu-ararat-yerevan/OGMRfhZo/"

Write-Host "Added row 82 to sheet '$($ws.Name)'"
